$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style/formatting from the existing last header cell (AB1)
# onto the new header cells so they match the bold/centered/bordered look.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the season record (Wins / Losses / Ties) for every data row.
for ($r = 2; $r -le 41; $r++) {
    $ws.Range("AC$r").Value = 87
    $ws.Range("AD$r").Value = 75
    $ws.Range("AE$r").Value = 0
}
